$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Save" header column in H1 (after existing G1 "sum" column)
$ws.Range("H1").Value = "Save"

# Match H1's formatting to the other header cells (bold, bordered, centered)
# by copying G1's format rather than constructing a brand-new style entry.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)

# Add corresponding data value in H2
$ws.Range("H2").Value = 0
